$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 10000
$ws.Range("M7").Value = -9888
# Row 14
$ws.Range("H14").Value = 10000
$ws.Range("I14").Value = 10000
$ws.Range("K14").Value = 10000
$ws.Range("M14").Value = -9809
# Row 19
$ws.Range("H19").Value = 1173.4
$ws.Range("I19").Value = 674.5
$ws.Range("J19").Value = 1250.1538
$ws.Range("K19").Value = 674.5
$ws.Range("L19").Value = 1250.1538
$ws.Range("M19").Value = -499.5
$ws.Range("N19").Value = -1600.1538
# Row 32
$ws.Range("H32").Value = 5362126
$ws.Range("I32").Value = 988.5
$ws.Range("J32").Value = 6336878
$ws.Range("K32").Value = 988.5
$ws.Range("L32").Value = 6336878
$ws.Range("M32").Value = -662.5
$ws.Range("N32").Value = -6337530
# Row 116
$ws.Range("H116").Value = 3693.647
$ws.Range("I116").Value = 4001.0908
$ws.Range("J116").Value = 3130
$ws.Range("K116").Value = 4001.0908
$ws.Range("L116").Value = 3130
$ws.Range("M116").Value = -559.0907999999999
$ws.Range("N116").Value = -10014
# Row 132
$ws.Range("H132").Value = 1442.738
$ws.Range("I132").Value = 1013.6301
$ws.Range("J132").Value = 4290.4546
$ws.Range("K132").Value = 3040.8903
$ws.Range("L132").Value = 12871.3638
$ws.Range("M132").Value = -510.8903
$ws.Range("N132").Value = -17931.3638
# Row 137
$ws.Range("H137").Value = 6494.241
$ws.Range("I137").Value = 9085.723
$ws.Range("J137").Value = 2253.6365
$ws.Range("K137").Value = 27257.169
$ws.Range("L137").Value = 6760.9095
$ws.Range("M137").Value = -24707.169
$ws.Range("N137").Value = -11860.9095

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4425.8887
$ws.Range("I86").Value = 5522.3335
$ws.Range("J86").Value = 2233
$ws.Range("K86").Value = 5522.3335
$ws.Range("L86").Value = 2233
$ws.Range("M86").Value = -4399.3335
$ws.Range("N86").Value = -4479
# Row 89
$ws.Range("H89").Value = 4425.8887
$ws.Range("I89").Value = 5522.3335
$ws.Range("J89").Value = 2233
$ws.Range("K89").Value = 27611.6675
$ws.Range("L89").Value = 11165
$ws.Range("M89").Value = -21995.6675
$ws.Range("N89").Value = -22397
# Row 94
$ws.Range("H94").Value = 1192.4667
$ws.Range("I94").Value = 837.4
$ws.Range("J94").Value = 1370
$ws.Range("K94").Value = 837.4
$ws.Range("L94").Value = 1370
$ws.Range("M94").Value = -386.4
$ws.Range("N94").Value = -2272
# Row 99
$ws.Range("H99").Value = 5015.385
$ws.Range("I99").Value = 5015.385
$ws.Range("K99").Value = 5015.385
$ws.Range("M99").Value = -3517.385
# Row 105
$ws.Range("H105").Value = 3126583.2
$ws.Range("I105").Value = 1513.1072
$ws.Range("K105").Value = 1513.1072
$ws.Range("M105").Value = 233.8928000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 372.2857
$ws.Range("I10").Value = 301.4
$ws.Range("J10").Value = 549.5
$ws.Range("K10").Value = 301.4
$ws.Range("L10").Value = 549.5
$ws.Range("M10").Value = -162.4
$ws.Range("N10").Value = -827.5
# Row 31
$ws.Range("H31").Value = 2256.4443
$ws.Range("I31").Value = 1052.9354
$ws.Range("J31").Value = 4921.357
$ws.Range("K31").Value = 1052.9354
$ws.Range("L31").Value = 4921.357
$ws.Range("M31").Value = -757.9354000000001
$ws.Range("N31").Value = -5511.357
# Row 34
$ws.Range("H34").Value = 2256.4443
$ws.Range("I34").Value = 1052.9354
$ws.Range("J34").Value = 4921.357
$ws.Range("K34").Value = 1052.9354
$ws.Range("L34").Value = 4921.357
$ws.Range("M34").Value = -850.9354000000001
$ws.Range("N34").Value = -5325.357
# Row 58
$ws.Range("H58").Value = 4382.6313
$ws.Range("I58").Value = 8593.308000000001
$ws.Range("J58").Value = 2193.08
$ws.Range("K58").Value = 8593.308000000001
$ws.Range("L58").Value = 2193.08
$ws.Range("M58").Value = -8390.308000000001
$ws.Range("N58").Value = -2599.08
# Row 62
$ws.Range("H62").Value = 3476039.5
$ws.Range("J62").Value = 4267.4614
$ws.Range("L62").Value = 4267.4614
$ws.Range("N62").Value = -5515.4614
# Row 65
$ws.Range("H65").Value = 3476039.5
$ws.Range("J65").Value = 4267.4614
$ws.Range("L65").Value = 21337.307
$ws.Range("N65").Value = -27577.307
# Row 94
$ws.Range("H94").Value = 5061.4287
$ws.Range("I94").Value = 1103.6
$ws.Range("J94").Value = 8659.454
$ws.Range("K94").Value = 1103.6
$ws.Range("L94").Value = 8659.454
$ws.Range("M94").Value = -652.5999999999999
$ws.Range("N94").Value = -9561.454
# Row 99
$ws.Range("H99").Value = 127082
$ws.Range("I99").Value = 201522.4
$ws.Range("J99").Value = 3014.6667
$ws.Range("K99").Value = 201522.4
$ws.Range("L99").Value = 3014.6667
$ws.Range("M99").Value = -200024.4
$ws.Range("N99").Value = -6010.6667
# Row 126
$ws.Range("H126").Value = 127082
$ws.Range("I126").Value = 201522.4
$ws.Range("J126").Value = 3014.6667
$ws.Range("K126").Value = 604567.2
$ws.Range("L126").Value = 9044.000100000001
$ws.Range("M126").Value = -602097.2
$ws.Range("N126").Value = -13984.0001
# Row 132
$ws.Range("H132").Value = 2127.0312
$ws.Range("I132").Value = 914.65
$ws.Range("K132").Value = 2743.95
$ws.Range("M132").Value = -213.9499999999998
# Row 136
$ws.Range("H136").Value = 4382.6313
$ws.Range("I136").Value = 8593.308000000001
$ws.Range("J136").Value = 2193.08
$ws.Range("K136").Value = 25779.924
$ws.Range("L136").Value = 6579.24
$ws.Range("M136").Value = -23229.924
$ws.Range("N136").Value = -11679.24

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1035.5862
$ws.Range("I5").Value = 390.8125
$ws.Range("J5").Value = 1829.1538
$ws.Range("K5").Value = 1172.4375
$ws.Range("L5").Value = 5487.4614
$ws.Range("M5").Value = -1060.4375
$ws.Range("N5").Value = -5711.4614
# Row 86
$ws.Range("H86").Value = 640.64703
$ws.Range("I86").Value = 676.6667
$ws.Range("J86").Value = 370.5
$ws.Range("K86").Value = 2030.0001
$ws.Range("L86").Value = 1111.5
$ws.Range("M86").Value = -844.0001
$ws.Range("N86").Value = -3483.5
# Row 89
$ws.Range("H89").Value = 640.64703
$ws.Range("I89").Value = 676.6667
$ws.Range("J89").Value = 370.5
$ws.Range("K89").Value = 6090.0003
$ws.Range("L89").Value = 3334.5
$ws.Range("M89").Value = -162.0002999999997
$ws.Range("N89").Value = -15190.5
# Row 107
$ws.Range("H107").Value = 665.9091
$ws.Range("I107").Value = 264.14285
$ws.Range("J107").Value = 853.4
$ws.Range("K107").Value = 792.4285500000001
$ws.Range("L107").Value = 2560.2
$ws.Range("M107").Value = 1127.57145
$ws.Range("N107").Value = -6400.2
# Row 113
$ws.Range("H113").Value = 9804446
$ws.Range("I113").Value = 19231314
$ws.Range("J113").Value = 502.76
$ws.Range("K113").Value = 57693942
$ws.Range("L113").Value = 1508.28
$ws.Range("M113").Value = -57691772
$ws.Range("N113").Value = -5848.28
# Row 131
$ws.Range("H131").Value = 1190.3226
$ws.Range("I131").Value = 480
$ws.Range("J131").Value = 1280.7273
$ws.Range("K131").Value = 1440
$ws.Range("L131").Value = 3842.1819
$ws.Range("M131").Value = 3600
$ws.Range("N131").Value = -13922.1819
# Row 132
$ws.Range("H132").Value = 3802.0938
$ws.Range("I132").Value = 2175.1333
$ws.Range("J132").Value = 5237.647
$ws.Range("K132").Value = 19576.1997
$ws.Range("L132").Value = 47138.823
$ws.Range("M132").Value = -17046.1997
$ws.Range("N132").Value = -52198.823
# Row 135
$ws.Range("H135").Value = 1035.5862
$ws.Range("I135").Value = 390.8125
$ws.Range("J135").Value = 1829.1538
$ws.Range("K135").Value = 3517.3125
$ws.Range("L135").Value = 16462.3842
$ws.Range("M135").Value = -982.3125
$ws.Range("N135").Value = -21532.3842

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 4000
$ws.Range("J3").Value = 4000
$ws.Range("L3").Value = 4000
$ws.Range("N3").Value = -4224
# Row 15
$ws.Range("H15").Value = 4000
$ws.Range("J15").Value = 4000
$ws.Range("L15").Value = 4000
$ws.Range("N15").Value = -4340
# Row 132
$ws.Range("H132").Value = 10255.379
$ws.Range("I132").Value = 4850.3
$ws.Range("J132").Value = 13100.158
$ws.Range("K132").Value = 14550.9
$ws.Range("L132").Value = 39300.474
$ws.Range("M132").Value = -12020.9
$ws.Range("N132").Value = -44360.474

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 775.8
$ws.Range("I7").Value = 775.8
$ws.Range("K7").Value = 775.8
$ws.Range("M7").Value = -662.8
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
# Row 9
$ws.Range("H9").Value = 530
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 530
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 530
$ws.Range("N9").Value = -810
$ws.Range("M9").ClearContents()
# Row 12
$ws.Range("H12").Value = 3166.6667
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 3500
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3500
$ws.Range("M12").Value = -2858
$ws.Range("N12").Value = -3784
# Row 81
$ws.Range("H81").Value = 2081
$ws.Range("I81").Value = 936.375
$ws.Range("K81").Value = 1872.75
$ws.Range("M81").Value = -811.75
# Row 84
$ws.Range("H84").Value = 2081
$ws.Range("I84").Value = 936.375
$ws.Range("K84").Value = 9363.75
$ws.Range("M84").Value = -4059.75
# Row 113
$ws.Range("H113").Value = 313.375
$ws.Range("I113").Value = 303.63635
$ws.Range("J113").Value = 334.8
$ws.Range("K113").Value = 910.90905
$ws.Range("L113").Value = 1004.4
$ws.Range("M113").Value = 1259.09095
$ws.Range("N113").Value = -5344.4
# Row 132
$ws.Range("H132").Value = 3683.923
$ws.Range("I132").Value = 2197
$ws.Range("K132").Value = 6591
$ws.Range("M132").Value = -4061
